$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation rows to append to the "Translations" sheet.
# Columns: A = Name, B = SourceLanguage, C = TargetLanguage, D = Status
$rows = @(
    @("Cat",             "Cat",             "Chat",              "Need review"),
    @("Password",        "Password",        "Mot de passe",      "Need review"),
    @("Run if you like",  "Run if you like", "Cours si tu veux",  "Need review"),
    @("House",           "House",           "Loger",             "Need review"),
    @("Wood",            "Wood",            "Bois",              "Need review")
)

$startRow = 3
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = $rows[$i][3]
}

$wb.Save()
